# Applies the "Updated symbol list" commit: refreshed price/volume/hour
# data for existing coin rows and swapped the BOLO / CoinbaseStockToken
# rows (with their own refreshed data) at the bottom of the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. All values must be written as
# literal text (matching the original inlineStr cells), so numeric-looking
# ones get a leading apostrophe plus ClearFormats() to strip the resulting
# quote-prefix style and keep the cell looking like a plain, unstyled cell.
$updates = [ordered]@{
    'D2' = '296.82'
    'E2' = '2.38%'
    'G2' = '21'
    'D3' = '41.29'
    'E3' = '2.23%'
    'G3' = '21'
    'D4' = '5.036'
    'E4' = '-0.19%'
    'G4' = '21'
    'D5' = '0.07540'
    'E5' = '3.32%'
    'G5' = '21'
    'D6' = '1.588'
    'E6' = '1.47%'
    'G6' = '21'
    'E7' = '0.84%'
    'G7' = '21'
    'D8' = '2.423'
    'E8' = '1.09%'
    'G8' = '21'
    'D9' = '0.1212'
    'E9' = '4.89%'
    'G9' = '21'
    'D10' = '0.1825'
    'E10' = '5.63%'
    'G10' = '21'
    'D11' = '0.08933'
    'E11' = '3.50%'
    'G11' = '21'
    'D12' = '0.04025'
    'E12' = '-3.52%'
    'G12' = '21'
    'E13' = '0.25%'
    'G13' = '21'
    'D14' = '0.001287'
    'E14' = '1.34%'
    'G14' = '21'
    'D15' = '0.005961'
    'E15' = '2.97%'
    'G15' = '21'
    'D16' = '3.340'
    'E16' = '-1.61%'
    'G16' = '21'
    'D17' = '4.377'
    'E17' = '2.29%'
    'G17' = '21'
    'D18' = '0.3318'
    'E18' = '1.25%'
    'G18' = '21'
    'D19' = '7.982'
    'E19' = '1.91%'
    'G19' = '21'
    'D20' = '0.1418'
    'E20' = '4.95%'
    'G20' = '21'
    'D21' = '0.3004'
    'E21' = '4.19%'
    'G21' = '21'
    'D22' = '0.04051'
    'E22' = '4.76%'
    'G22' = '21'
    'D23' = '0.001267'
    'E23' = '-0.13%'
    'G23' = '21'
    'D24' = '0.004025'
    'E24' = '3.89%'
    'G24' = '21'
    'D25' = '0.0001232'
    'E25' = '-3.83%'
    'G25' = '21'
    'E26' = '-0.03%'
    'G26' = '21'
    'G27' = '21'
    'G28' = '21'
    'G29' = '21'
    'G30' = '21'
    'G31' = '21'
    'G32' = '21'
    'G33' = '21'
    'G34' = '21'
    'G35' = '21'
    'G36' = '21'
    'G37' = '21'
    'D38' = '0.02415'
    'E38' = '4.05%'
    'G38' = '21'
    'D39' = '0.05211'
    'E39' = '4.89%'
    'G39' = '21'
    'D40' = '0.006550'
    'E40' = '-3.75%'
    'G40' = '21'
    'D41' = '0.007781'
    'E41' = '1.00%'
    'G41' = '21'
    'D42' = '0.1329'
    'E42' = '4.54%'
    'G42' = '21'
    'D43' = '0.007577'
    'E43' = '3.03%'
    'G43' = '21'
    'D44' = '0.007861'
    'E44' = '11.25%'
    'G44' = '21'
    'D45' = '0.3213'
    'E45' = '11.13%'
    'G45' = '21'
    'D46' = '0.00006787'
    'E46' = '5.70%'
    'G46' = '21'
    'D47' = '0.00000000751'
    'E47' = '0.07%'
    'G47' = '21'
    'B48' = 'BOLO'
    'C48' = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
    'D48' = '0.04616'
    'E48' = '89.08%'
    'G48' = '21'
    'B49' = 'CoinbaseStockToken'
    'C49' = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
    'D49' = '0.004206'
    'E49' = '0.11%'
    'G49' = '21'
    'D50' = '0.00002103'
    'E50' = '0.07%'
    'G50' = '21'
    'D51' = '0.0002003'
    'E51' = '0.07%'
    'G51' = '21'
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)
    if ($value -match '^-?\d+(\.\d+)?%?$') {
        # Numeric-looking text: force literal text via quote-prefix, then
        # clear the format Excel auto-applied so no style is left behind.
        $cell.Value = "'" + $value
        $cell.ClearFormats()
    } else {
        $cell.Value = $value
    }
}
